$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 76, shifting existing rows 76-87 down to 77-88.
$ws.Rows.Item(76).Insert()

# Populate the new row 76 with a copy of the (now shifted) row-77 context
# values, plus the new date/price figures from this week's entry.
$ws.Range("A76").Value = 2
$ws.Range("B76").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C76").Value = "Coquimbo"
$ws.Range("D76").Value = 45077
$ws.Range("E76").Value = 4
$ws.Range("F76").Value = 100112026
$ws.Range("G76").Value = "Haba"
$ws.Range("H76").Value = "Sin especificar"
$ws.Range("I76").Value = "Primera"
$ws.Range("J76").Value = 700
$ws.Range("K76").Value = 12000
$ws.Range("L76").Value = 14000
$ws.Range("M76").Value = 13000
$ws.Range("N76").Value = "$/saco 25 kilos"
$ws.Range("O76").Value = "Provincia de Limarí"
$ws.Range("P76").Value = 520
$ws.Range("Q76").Value = 25
$ws.Range("R76").Value = "Hortaliza"
